$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.937933333333332
$ws.Range("H2").Value = 26.8138
$ws.Range("I2").Value = 0.2302024600837126
$ws.Range("J2").Value = 0.2302024600837126
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 344.5812538682221
$ws.Range("R2").Value = 3101.231284814
$ws.Range("S2").Value = 0.1325540692832748
$ws.Range("T2").Value = 0.1325540692832747
$ws.Range("G3").Value = 8.937933333333332
$ws.Range("H3").Value = 26.8138
$ws.Range("I3").Value = 0.2302024600837126
$ws.Range("J3").Value = 0.2302024600837126
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 48.00274404293333
$ws.Range("R3").Value = 432.0246963864
$ws.Range("S3").Value = 0.0184657725521182
$ws.Range("T3").Value = 0.0184657725521182
$ws.Range("G4").Value = 8.937933333333332
$ws.Range("H4").Value = 26.8138
$ws.Range("I4").Value = 0.2302024600837126
$ws.Range("J4").Value = 0.2302024600837126
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 205.8393682525555
$ws.Range("R4").Value = 1852.554314273
$ws.Range("S4").Value = 0.07918261824831964
$ws.Range("T4").Value = 0.07918261824831964
$ws.Range("I5").Value = 0.5278886986241245
$ws.Range("J5").Value = 0.5278886986241244
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 790.1763934608566
$ws.Range("R5").Value = 7111.587541147709
$ws.Range("S5").Value = 0.3039663221054811
$ws.Range("T5").Value = 0.303966322105481
$ws.Range("I6").Value = 0.5278886986241245
$ws.Range("J6").Value = 0.5278886986241244
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("R6").Value = 990.697296049596
$ws.Range("S6").Value = 0.04234478049488246
$ws.Range("T6").Value = 0.04234478049488246
$ws.Range("I7").Value = 0.5278886986241245
$ws.Range("J7").Value = 0.5278886986241244
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 472.0204822873716
$ws.Range("R7").Value = 4248.184340586345
$ws.Range("S7").Value = 0.1815775960237609
$ws.Range("T7").Value = 0.1815775960237609
$ws.Range("G8").Value = 9.392449999999998
$ws.Range("H8").Value = 28.17735
$ws.Range("I8").Value = 0.241908841292163
$ws.Range("J8").Value = 0.2419088412921629
$ws.Range("M8").Value = 38.55267666666666
$ws.Range("N8").Value = 115.65803
$ws.Range("O8").Value = 0.5758151725879548
$ws.Range("P8").Value = 0.5758151725879548
$ws.Range("Q8").Value = 362.1040879578333
$ws.Range("R8").Value = 3258.9367916205
$ws.Range("S8").Value = 0.139294781199199
$ws.Range("T8").Value = 0.139294781199199
$ws.Range("G9").Value = 9.392449999999998
$ws.Range("H9").Value = 28.17735
$ws.Range("I9").Value = 0.241908841292163
$ws.Range("J9").Value = 0.2419088412921629
$ws.Range("O9").Value = 0.08021535714867321
$ws.Range("P9").Value = 0.08021535714867323
$ws.Range("Q9").Value = 50.4438057962
$ws.Range("R9").Value = 453.9942521658
$ws.Range("S9").Value = 0.01940480410167256
$ws.Range("T9").Value = 0.01940480410167256
$ws.Range("G10").Value = 9.392449999999998
$ws.Range("H10").Value = 28.17735
$ws.Range("I10").Value = 0.241908841292163
$ws.Range("J10").Value = 0.2419088412921629
$ws.Range("M10").Value = 23.02986166666667
$ws.Range("N10").Value = 69.089585
$ws.Range("O10").Value = 0.3439694702633719
$ws.Range("P10").Value = 0.3439694702633719
$ws.Range("Q10").Value = 216.3068242110833
$ws.Range("R10").Value = 1946.76141789975
$ws.Range("S10").Value = 0.08320925599129141
$ws.Range("T10").Value = 0.08320925599129141
